$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.768.58"
$ws.Range("E2").Value = "  +9.44%  "
$ws.Range("D3").Value = "1.780.15"
$ws.Range("E3").Value = "  +6.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.553"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.278"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0656"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0924"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").Value = "2.041.17"
$ws.Range("E13").Value = "  +6.28%  "
$ws.Range("D14").Value = "1.778.08"
$ws.Range("E14").Value = "  +6.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.626"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.01%  "
$ws.Range("D16").Value = "33.653.34"
$ws.Range("E16").Value = "  +9.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "250.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "
$ws.Range("D21").Value = "0.0₃0735"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0508"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("E33").Value = "  +4.54%  "
$ws.Range("E34").Value = "  +6.07%  "
$ws.Range("D35").Value = "1.485.39"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.617"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0185"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.879"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.89%  "
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0514"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("D47").Value = "1.931.16"
$ws.Range("E47").Value = "  +6.40%  "
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "
